$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that originally sat
#    right after the "Play Free Five Pirates Slot Game Review" heading.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Meta description:*") {
        $para.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2. Insert a new paragraph with bold text "Play Free Five Pirates
#    Slot Game Review" right before the closing "Prompt: ..." paragraph.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r/>
            <w:r>
              <w:rPr>
                <w:b/>
              </w:rPr>
              <w:t>Play Free Five Pirates Slot Game Review</w:t>
            </w:r>
          </w:p>
          <w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$insertPoint.InsertXML($xml)

# The insertion above leaves behind one stray empty paragraph right
# after the new bold paragraph (an artifact of splitting the
# paragraph); remove it so the structure matches exactly.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if (($p.Range.End - $p.Range.Start) -le 1) {
        if ($i -lt $d.Paragraphs.Count -and $d.Paragraphs.Item($i + 1).Range.Text -like "Prompt:*") {
            $p.Range.Delete()
            break
        }
    }
}

# ------------------------------------------------------------------
# 3. Replace the old "Prompt: ..." image-generation text with the new
#    meta-description copy, keeping the italic run formatting intact.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Prompt: Create an eye-catching feature image for the online slot game, " + [char]34 + "Five Pirates" + [char]34 + ", which showcases a happy Maya warrior with glasses in a cartoon style. The image should be bright and colorful, with the Maya warrior as the main focus. The warrior should exude a sense of adventure, bravery, and fun. Be sure to include elements of the pirate theme, such as treasure maps, chests, parrots, cannons, and gold doubloons in the background or surrounding the Maya warrior. Use creative liberty to make the image exciting and attention-grabbing, while still maintaining the integrity of the game's theme.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Explore the ocean and play Five Pirates slot game from Lightning Box Games. Featuring 1,024 ways to win and bonuses, including free spins. Play for free now!",
    2)
